$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 242, shifting existing rows 242:320 down to 244:322
$ws.Rows("242:243").Insert()

# Copy the date number format (style) used by the rest of column D onto the new D cells
$ws.Range("D244").Copy()
$ws.Range("D242:D243").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 242 values
$ws.Range("A242").Value = 1
$ws.Range("B242").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C242").Value = "Arica y Parinacota"
$ws.Range("D242").Value = 44722
$ws.Range("E242").Value = 15
$ws.Range("F242").Value = 100112043
$ws.Range("G242").Value = "Pepino ensalada"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 120
$ws.Range("K242").Value = 15000
$ws.Range("L242").Value = 16000
$ws.Range("M242").Value = 15500
$ws.Range("N242").Value = "$/caja 70 unidades"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 221
$ws.Range("Q242").Value = 70
$ws.Range("R242").Value = "Hortaliza"

# Row 243 values
$ws.Range("A243").Value = 1
$ws.Range("B243").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C243").Value = "Arica y Parinacota"
$ws.Range("D243").Value = 44722
$ws.Range("E243").Value = 15
$ws.Range("F243").Value = 100112043
$ws.Range("G243").Value = "Pepino ensalada"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Segunda"
$ws.Range("J243").Value = 150
$ws.Range("K243").Value = 13000
$ws.Range("L243").Value = 14000
$ws.Range("M243").Value = 13500
$ws.Range("N243").Value = "$/caja 100 unidades"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 135
$ws.Range("Q243").Value = 100
$ws.Range("R243").Value = "Hortaliza"
